# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to H/I/J/K/L/M/N columns across several sheets/rows
# as described by the source diff (Hyperion_Profits market-data recalculation).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 884.0816
$ws.Range("I15").Value = 884.0816
$ws.Range("K15").Value = 2652.2448
$ws.Range("M15").Value = -2483.2448

$ws.Range("H100").Value = 857.1429000000001
$ws.Range("I100").Value = 833.3333
$ws.Range("K100").Value = 833.3333
$ws.Range("M100").Value = -292.3333

$ws.Range("H107").Value = 55556280
$ws.Range("I107").Value = 66667336
$ws.Range("K107").Value = 66667336
$ws.Range("M107").Value = -66665416

$ws.Range("H111").Value = 9261738
$ws.Range("J111").Value = 1207
$ws.Range("L111").Value = 3621
$ws.Range("N111").Value = -9755

$ws.Range("H133").Value = 124330.43
$ws.Range("J133").Value = 124330.43
$ws.Range("L133").Value = 124330.43
$ws.Range("N133").Value = -134450.43

$ws.Range("H140").Value = 96481.71000000001
$ws.Range("J140").Value = 96481.71000000001
$ws.Range("L140").Value = 96481.71000000001
$ws.Range("N140").Value = -106841.71

$ws.Range("H141").Value = 18686.125
$ws.Range("J141").Value = 4999.5
$ws.Range("L141").Value = 14998.5
$ws.Range("N141").Value = -25358.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10059.228
$ws.Range("I32").Value = 6176.122
$ws.Range("K32").Value = 6176.122
$ws.Range("M32").Value = -5889.122

$ws.Range("H132").Value = 2932.423
$ws.Range("I132").Value = 2302.25
$ws.Range("K132").Value = 6906.75
$ws.Range("M132").Value = -4376.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3777.611
$ws.Range("I20").Value = 3104.5
$ws.Range("K20").Value = 3104.5
$ws.Range("M20").Value = -2857.5

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H134").Value = 3610.919
$ws.Range("I134").Value = 1405.8422
$ws.Range("K134").Value = 4217.5266
$ws.Range("M134").Value = -1682.5266

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1611.6428
$ws.Range("I16").Value = 961.1818
$ws.Range("J16").Value = 3996.6667
$ws.Range("K16").Value = 961.1818
$ws.Range("L16").Value = 3996.6667
$ws.Range("M16").Value = -674.1818
$ws.Range("N16").Value = -4570.6667

$ws.Range("H22").Value = 799
$ws.Range("I22").Value = 799
$ws.Range("K22").Value = 799
$ws.Range("M22").Value = -449

$ws.Range("H113").Value = 1611.6428
$ws.Range("I113").Value = 961.1818
$ws.Range("J113").Value = 3996.6667
$ws.Range("K113").Value = 961.1818
$ws.Range("L113").Value = 3996.6667
$ws.Range("M113").Value = 1208.8182
$ws.Range("N113").Value = -8336.6667

$ws.Range("H132").Value = 78020.56
$ws.Range("J132").Value = 91768.3
$ws.Range("L132").Value = 275304.9
$ws.Range("N132").Value = -280364.9

$ws.Range("H141").Value = 180271.67
$ws.Range("J141").Value = 180271.67
$ws.Range("L141").Value = 180271.67
$ws.Range("N141").Value = -190631.67

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13839648
$ws.Range("I4").Value = 15096658
$ws.Range("K4").Value = 45289974
$ws.Range("M4").Value = -45289862

$ws.Range("H86").Value = 200
$ws.Range("J86").Value = 300
$ws.Range("L86").Value = 900
$ws.Range("N86").Value = -3272

$ws.Range("H89").Value = 200
$ws.Range("J89").Value = 300
$ws.Range("L89").Value = 2700
$ws.Range("N89").Value = -14556

$ws.Range("H107").Value = 2500
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws.Range("H110").Value = 19900.1
$ws.Range("I110").Value = 1001
$ws.Range("K110").Value = 3003
$ws.Range("M110").Value = 1087

$ws.Range("H115").Value = 99
$ws.Range("I115").Value = 99
$ws.Range("K115").Value = 297
$ws.Range("M115").Value = 878

$ws.Range("H131").Value = 11577101
$ws.Range("I131").Value = 5557333.5
$ws.Range("K131").Value = 16672000.5
$ws.Range("M131").Value = -16666960.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13343573
$ws.Range("I70").Value = 15395585
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 15395585
$ws.Range("L70").Value = 5500
$ws.Range("M70").Value = -15395315
$ws.Range("N70").Value = -6040

$ws.Range("H73").Value = 13343573
$ws.Range("I73").Value = 15395585
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 15395585
$ws.Range("L73").Value = 5500
$ws.Range("M73").Value = -15394649
$ws.Range("N73").Value = -7372

$ws.Range("H92").Value = 22309.6
$ws.Range("J92").Value = 22309.6
$ws.Range("L92").Value = 22309.6
$ws.Range("N92").Value = -26053.6

$ws.Range("H95").Value = 47449.5
$ws.Range("J95").Value = 47449.5
$ws.Range("L95").Value = 47449.5
$ws.Range("N95").Value = -52941.5

$ws.Range("H102").Value = 3778287.5
$ws.Range("I102").Value = 7408954
$ws.Range("K102").Value = 7408954
$ws.Range("M102").Value = -7407332

$ws.Range("H126").Value = 5858239.5
$ws.Range("I126").Value = 9093589
$ws.Range("J126").Value = 4906666.5
$ws.Range("K126").Value = 27280767
$ws.Range("L126").Value = 14719999.5
$ws.Range("M126").Value = -27278297
$ws.Range("N126").Value = -14724939.5

$ws.Range("H132").Value = 2948
$ws.Range("I132").Value = 2879.7666
$ws.Range("K132").Value = 8639.299800000001
$ws.Range("M132").Value = -6109.299800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 50397.445
$ws.Range("I22").Value = 222595.5
$ws.Range("J22").Value = 1198
$ws.Range("K22").Value = 222595.5
$ws.Range("L22").Value = 1198
$ws.Range("M22").Value = -222300.5
$ws.Range("N22").Value = -1788

$ws.Range("H27").Value = 50397.445
$ws.Range("I27").Value = 222595.5
$ws.Range("J27").Value = 1198
$ws.Range("K27").Value = 222595.5
$ws.Range("L27").Value = 1198
$ws.Range("M27").Value = -222488.5
$ws.Range("N27").Value = -1412

$ws.Range("H46").Value = 4427.0454
$ws.Range("I46").Value = 704.375
$ws.Range("J46").Value = 6554.2856
$ws.Range("K46").Value = 704.375
$ws.Range("L46").Value = 6554.2856
$ws.Range("M46").Value = -516.375
$ws.Range("N46").Value = -6930.2856

$ws.Range("H82").Value = 4276492.5
$ws.Range("I82").Value = 5558381
$ws.Range("K82").Value = 5558381
$ws.Range("M82").Value = -5558020

$ws.Range("H85").Value = 4276492.5
$ws.Range("I85").Value = 5558381
$ws.Range("K85").Value = 5558381
$ws.Range("M85").Value = -5557133

$ws.Range("H132").Value = 4101.6855
$ws.Range("I132").Value = 3354.32
$ws.Range("J132").Value = 5970.1
$ws.Range("K132").Value = 10062.96
$ws.Range("L132").Value = 17910.3
$ws.Range("M132").Value = -7532.960000000001
$ws.Range("N132").Value = -22970.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 16332
$ws.Range("I55").Value = 14010.5
$ws.Range("J55").Value = 20975
$ws.Range("K55").Value = 14010.5
$ws.Range("L55").Value = 20975
$ws.Range("M55").Value = -13733.5
$ws.Range("N55").Value = -21529

$ws.Range("H68").Value = 23000
$ws.Range("J68").Value = 23000
$ws.Range("L68").Value = 23000
$ws.Range("N68").Value = -24622

$ws.Range("H71").Value = 23000
$ws.Range("J71").Value = 23000
$ws.Range("L71").Value = 69000
$ws.Range("N71").Value = -77112

$ws.Range("H107").Value = 41667828
$ws.Range("I107").Value = 55555944
$ws.Range("J107").Value = 3480.6667
$ws.Range("K107").Value = 166667832
$ws.Range("L107").Value = 10442.0001
$ws.Range("M107").Value = -166665912
$ws.Range("N107").Value = -14282.0001

$ws.Range("H122").Value = 4077.2
$ws.Range("J122").Value = 4097.8335
$ws.Range("L122").Value = 12293.5005
$ws.Range("N122").Value = -17193.5005

$ws.Range("H132").Value = 21509964

$ws.Range("H136").Value = 1993.0682
$ws.Range("I136").Value = 1450.0555
$ws.Range("J136").Value = 4436.625
$ws.Range("K136").Value = 4350.166499999999
$ws.Range("L136").Value = 13309.875
$ws.Range("M136").Value = -1800.166499999999
$ws.Range("N136").Value = -18409.875

